$wb = $excel.ActiveWorkbook

# "전국" (nationwide) sheet is the first sheet (index 1)
$wsNational = $wb.Worksheets.Item(1)

# Reset all result cells (B2:Q19) on the nationwide sheet to 0, then set
# the special "합계" (total) row 2 figures that remained non-zero.
$wsNational.Range("B2:Q19").Value2 = 0

$wsNational.Range("D2").Value2 = 7000
$wsNational.Range("E2").Value2 = 3000
$wsNational.Range("P2").Value2 = 10000

# Make the nationwide sheet the active tab/sheet again, with E3 selected.
# (Seoul's own selection, T17, is untouched - it just stops being the
# active tab once another sheet is activated/selected.)
[void]$wsNational.Activate()
[void]$wsNational.Range("E3").Select()
